$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 29.75
$data[0,1] = 29.64127922058105
$data[0,2] = -0.1087207794189453
$data[0,3] = 0.01182020787746296
$data[1,0] = 29.84
$data[1,1] = 29.90331840515137
$data[1,2] = 0.06331840515136378
$data[1,3] = 0.004009220430912251
$data[2,0] = 29.81
$data[2,1] = 29.90265464782715
$data[2,2] = 0.09265464782714616
$data[2,3] = 0.008584883763972481
$data[3,0] = 29.92
$data[3,1] = 29.92182731628418
$data[3,2] = 0.001827316284177982
$data[3,3] = [double]"3.339084802422028e-06"
$data[4,0] = 29.98
$data[4,1] = 29.99698448181152
$data[4,2] = 0.01698448181151946
$data[4,3] = 0.0002884726224058353
$data[5,0] = 30.03999999999999
$data[5,1] = 30.15762138366699
$data[5,2] = 0.1176213836670001
$data[5,3] = 0.01383478989573965
$data[6,0] = 30.21000000000001
$data[6,1] = 30.14980888366699
$data[6,2] = -0.06019111633301577
$data[6,3] = 0.003622970485414638
$data[7,0] = 30.22
$data[7,1] = 30.26547050476074
$data[7,2] = 0.04547050476074332
$data[7,3] = 0.002067566803196782
$data[8,0] = 30.38
$data[8,1] = 30.23031425476074
$data[8,2] = -0.1496857452392533
$data[8,3] = 0.02240582232783063
$data[9,0] = 30.44
$data[9,1] = 30.44497489929199
$data[9,2] = 0.004974899291994461
$data[9,3] = [double]"2.474962296548699e-05"
$data[10,0] = 30.48
$data[10,1] = 30.40680503845215
$data[10,2] = -0.07319496154785554
$data[10,3] = 0.005357502395992051
$data[11,0] = 30.69
$data[11,1] = 30.44501304626465
$data[11,2] = -0.2449869537353493
$data[11,3] = 0.06001860750052617
$data[12,0] = 30.75
$data[12,1] = 30.42711448669434
$data[12,2] = -0.3228855133056641
$data[12,3] = 0.1042550547026622
$data[13,0] = 30.94
$data[13,1] = 30.66121482849121
$data[13,2] = -0.2787851715087868
$data[13,3] = 0.07772117185318367
$data[14,0] = 30.95
$data[14,1] = 30.80583000183105
$data[14,2] = -0.1441699981689482
$data[14,3] = 0.02078498837203452
$data[15,0] = 31.02
$data[15,1] = 31.21635246276855
$data[15,2] = 0.1963524627685587
$data[15,3] = 0.03855428963527822
$data[16,0] = 31.12
$data[16,1] = 31.35795402526855
$data[16,2] = 0.2379540252685501
$data[16,3] = 0.0566221181415058
$data[17,0] = 31.28
$data[17,1] = 31.45592308044434
$data[17,2] = 0.1759230804443348
$data[17,3] = 0.03094893023302389
$data[18,0] = 31.38
$data[18,1] = 31.29449272155762
$data[18,2] = -0.08550727844237827
$data[18,3] = 0.007311494666622407
$data[19,0] = 31.58
$data[19,1] = 31.48099327087402
$data[19,2] = -0.09900672912597486
$data[19,3] = 0.009802332412224158
$data[20,0] = 31.65000000000001
$data[20,1] = 32.05437088012695
$data[20,2] = 0.4043708801269474
$data[20,3] = 0.1635158086946421
$data[21,0] = 31.88
$data[21,1] = 32.60719680786133
$data[21,2] = 0.7271968078613327
$data[21,3] = 0.528815197363712
$data[22,0] = 32.28
$data[22,1] = 32.54343032836914
$data[22,2] = 0.2634303283691395
$data[22,3] = 0.06939553790467266
$data[23,0] = 32.45
$data[23,1] = 32.63210678100586
$data[23,2] = 0.1821067810058565
$data[23,3] = 0.03316287968831499
$data[24,0] = 32.84999999999999
$data[24,1] = 32.83679580688477
$data[24,2] = -0.01320419311522869
$data[24,3] = 0.0001743507158242527
$data[25,0] = 32.90000000000001
$data[25,1] = 33.05073165893555
$data[25,2] = 0.1507316589355412
$data[25,3] = 0.02272003300546031
$data[26,0] = 33.09999999999999
$data[26,1] = 32.78688430786133
$data[26,2] = -0.3131156921386662
$data[26,3] = 0.09804143666347598
$data[27,0] = 33.40000000000001
$data[27,1] = 33.60528182983398
$data[27,2] = 0.2052818298339787
$data[27,3] = 0.04214062965998659
$data[28,0] = 33.7
$data[28,1] = 33.57837295532227
$data[28,2] = -0.1216270446777372
$data[28,3] = 0.01479313799704029
$data[29,0] = 34.09999999999999
$data[29,1] = 33.80495834350586
$data[29,2] = -0.2950416564941349
$data[29,3] = 0.08704957906680312
$data[30,0] = 34.40000000000001
$data[30,1] = 34.3074836730957
$data[30,2] = -0.09251632690430256
$data[30,3] = 0.008559270743863777
$data[31,0] = 34.90000000000001
$data[31,1] = 35.01227951049805
$data[31,2] = 0.1122795104980412
$data[31,3] = 0.01260668847767974
$data[32,0] = 35.3
$data[32,1] = 35.66598129272461
$data[32,2] = 0.3659812927246122
$data[32,3] = 0.1339423066243783
$data[33,0] = 35.7
$data[33,1] = 35.98637771606445
$data[33,2] = 0.2863777160644503
$data[33,3] = 0.08201219625829091
$data[34,0] = 36.3
$data[34,1] = 35.96036148071289
$data[34,2] = -0.3396385192871065
$data[34,3] = 0.1153543237835382
$data[35,0] = 36.8
$data[35,1] = 36.60099411010742
$data[35,2] = -0.1990058898925753
$data[35,3] = 0.03960334421193579
$data[36,0] = 37.3
$data[36,1] = 37.27934646606445
$data[36,2] = -0.02065353393554403
$data[36,3] = 0.000426568464026669
$data[37,0] = 37.90000000000001
$data[37,1] = 38.10490798950195
$data[37,2] = 0.2049079895019474
$data[37,3] = 0.04198728416173021
$data[38,0] = 38.5
$data[38,1] = 38.57924270629883
$data[38,2] = 0.07924270629882812
$data[38,3] = 0.006279406501562335
$data[39,0] = 38.90000000000001
$data[39,1] = 39.20572280883789
$data[39,2] = 0.3057228088378849
$data[39,3] = 0.09346643584372594
$data[40,0] = 39.40000000000001
$data[40,1] = 39.69402694702148
$data[40,2] = 0.2940269470214787
$data[40,3] = 0.08645184557477144
$data[41,0] = 39.90000000000001
$data[41,1] = 39.77173233032227
$data[41,2] = -0.1282676696777401
$data[41,3] = 0.01645259508455784
$data[42,0] = 40.09999999999999
$data[42,1] = 39.85319900512695
$data[42,2] = -0.2468009948730412
$data[42,3] = 0.0609107310703229
$data[43,0] = 40.59999999999999
$data[43,1] = 40.23218154907227
$data[43,2] = -0.3678184509277287
$data[43,3] = 0.135290412842874
$data[44,0] = 40.90000000000001
$data[44,1] = 40.4200325012207
$data[44,2] = -0.4799674987793026
$data[44,3] = 0.2303687998844598
$data[45,0] = 41.20000000000001
$data[45,1] = 40.79880142211914
$data[45,2] = -0.4011985778808693
$data[45,3] = 0.160960298893632
$data[46,0] = 41.5
$data[46,1] = 40.99710464477539
$data[46,2] = -0.5028953552246094
$data[46,3] = 0.252903738306486
$data[47,0] = 41.8
$data[47,1] = 41.71234130859375
$data[47,2] = -0.08765869140624716
$data[47,3] = 0.00768404617905567
$data[48,0] = 42.2
$data[48,1] = 41.72422790527344
$data[48,2] = -0.4757720947265653
$data[48,3] = 0.2263590861205039
$data[49,0] = 42.7
$data[49,1] = 43.96305465698242
$data[49,2] = 1.263054656982419
$data[49,3] = 1.595307066524976

$ws.Range("B2:E51").Value = $data

$ws.Range("C52").Value = 0.1454766845702764
$ws.Range("E52").Value = 4.84477354914006
$ws.Range("E53").Value = 0.0968954709828012
